# Update "想去人数" (column F) figures, as generated by the site's data
# refresh at commit 456a3b4.
#
# The workbook has 4 sheets: 展览(1), 演出(2), 本地生活(3), 全部类型(4).
# The refreshed counts affect matching rows in the "展览" sheet and in the
# aggregate "全部类型" sheet (same events, different row numbers there).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (index 1) ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 261
$ws1.Range("F3").Value  = 1096
$ws1.Range("F5").Value  = 453
$ws1.Range("F6").Value  = 87
$ws1.Range("F7").Value  = 577
$ws1.Range("F9").Value  = 6894
$ws1.Range("F10").Value = 168
$ws1.Range("F15").Value = 1118
$ws1.Range("F16").Value = 16353
$ws1.Range("F18").Value = 1610
$ws1.Range("F20").Value = 340
$ws1.Range("F21").Value = 193
$ws1.Range("F22").Value = 120
$ws1.Range("F23").Value = 11452
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 1088
$ws1.Range("F26").Value = 4508
$ws1.Range("F27").Value = 365
$ws1.Range("F29").Value = 53
$ws1.Range("F30").Value = 852
$ws1.Range("F31").Value = 325

# ---- Sheet "全部类型" (index 4) ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 261
$ws4.Range("F3").Value  = 1096
$ws4.Range("F5").Value  = 453
$ws4.Range("F6").Value  = 87
$ws4.Range("F7").Value  = 577
$ws4.Range("F10").Value = 6894
$ws4.Range("F11").Value = 168
$ws4.Range("F17").Value = 1118
$ws4.Range("F18").Value = 16353
$ws4.Range("F20").Value = 1610
$ws4.Range("F22").Value = 340
$ws4.Range("F23").Value = 193
$ws4.Range("F24").Value = 120
$ws4.Range("F27").Value = 11452
$ws4.Range("F28").Value = 16
$ws4.Range("F29").Value = 1088
$ws4.Range("F30").Value = 4508
$ws4.Range("F31").Value = 365
$ws4.Range("F33").Value = 53
$ws4.Range("F34").Value = 852
$ws4.Range("F35").Value = 325
